$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# Row 7 ("Experimental") gained a value in column B: "true".
# Typing the literal text "true"/"false" straight into a cell makes Excel's
# COM layer coerce it to a real Boolean, which is not what the source
# document stores (it stores the literal text "true" as a shared string).
# Route it through a scratch cell forced to Text via a leading apostrophe,
# then copy/paste-special the *value* into B7 so the destination keeps its
# original (non quote-prefixed) style, and finally drop the scratch column.
$scratch = $ws1.Range("D1")
$scratch.Value = "'true"
$scratch.Copy()
$ws1.Range("B7").PasteSpecial(-4163)
$scratch.EntireColumn.Delete()

# Row 8 ("Date") value was refreshed to a newer timestamp string.
$ws1.Range("B8").Value = "2023-02-16T14:43:10-06:00"
